# Puskesmas 020 Pomalaa.xlsx - "add tabel format baru"
# Bumps the referenced table numbers (Tabel 4.2.3/4.2.4/4.2.5 -> 4.2.5/4.2.6/4.2.7)
# and the reporting year (2020 -> 2021) for the "Bab 4" sheet, and resets the
# saved view/scroll/selection state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table-number headings -------------------------------------------------
# H1: "Tabel 4.2.3" -> "Tabel 4.2.5"
$ws.Range("H1").Value = "Tabel 4.2.5"

# P1: rich text "Tabel" + " 4.2.4." -> "Tabel" + " 4.2.6."
$ws.Range("P1").Value = "Tabel 4.2.6."
$ws.Range("P1").Characters(6, 7).Font.Size = 9
$ws.Range("P1").Characters(6, 7).Font.Underline = $false

# W1: rich text "Tabel" + " 4.2.5." -> "Tabel" + " 4.2.7."
$ws.Range("W1").Value = "Tabel 4.2.7."
$ws.Range("W1").Characters(6, 7).Font.Size = 9
$ws.Range("W1").Characters(6, 7).Font.Underline = $false

# --- Year 2020 -> 2021 in the table titles ---------------------------------
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Pomalaa. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Pomalaa, 2021"
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Pomalaa, 2021"
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Pomalaa, 2021"

$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Pomalaa Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Pomalaa Subdistrict, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Pomalaa Subdistrict, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Pomalaa Subdistrict, 2021"

# --- Reset the saved view state (drop the stale scroll position/selection) -
$aw = $excel.ActiveWindow
$aw.Zoom = 85
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("A1").Select()
